$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark that used to sit right after
#    "Apologies:" - it gets relocated (see step 2 below).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Find the "Saving carts between sessions" bullet and add a brand new
#    top-level agenda item right after it.
$find = $d.Content
$find.Find.Execute("Saving carts between sessions", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$srcPara = $find.Paragraphs(1)

$tail = $srcPara.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()

$newPara = $srcPara.Next()
# Match the top-level bullet indent used by sibling agenda items (ilvl 0).
$newPara.Range.ListFormat.ListLevelNumber = 1

$boldLeadIn = "Discuss the Sprint 2 Stories and Any changes that have to be added before submission – HIGH"
$fullText = $boldLeadIn + "LY IMPORTANT "
$newPara.Range.Text = $fullText
$newPara.Range.Font.Bold = $true

# 3. Put the "_GoBack" bookmark back, this time between "...HIGH" and
#    "LY IMPORTANT " (i.e. where the cursor was left after the last edit).
$bmPos = $newPara.Range.Start + $boldLeadIn.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
